$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize columns D:M (user dragged/set column widths in the Rankings sheet)
$ws.Columns.Item(4).ColumnWidth = 10.67   # D -> 11.46484375
$ws.Columns.Item(5).ColumnWidth = 9.5     # E -> 10.3984375
$ws.Columns.Item(6).ColumnWidth = 10.33   # F -> 11.19921875 (style 6 preserved)
$ws.Columns.Item(7).ColumnWidth = 10.5    # G -> 11.3984375  (style 6 preserved)
$ws.Columns.Item(8).ColumnWidth = 10.33   # H -> 11.1328125
$ws.Columns.Item(9).ColumnWidth = 9.67    # I -> 10.53125
$ws.Columns.Item(10).ColumnWidth = 12.0   # J -> 12.86328125 (style 6 preserved)
$ws.Columns.Item(11).ColumnWidth = 12.5   # K -> 13.265625   (style 6 preserved)
$ws.Columns.Item(12).ColumnWidth = 11.5   # L -> 12.265625
$ws.Columns.Item(13).ColumnWidth = 0      # M -> 0

# Hide columns L and M
$ws.Columns.Item(12).Hidden = $true
$ws.Columns.Item(13).Hidden = $true

# Select column M (entire column), matching the new active selection
$ws.Columns.Item(13).Select()
